# Update column C (Förändrad) for rows 2-16 from 45175 to 45183.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 16; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45175) {
        $cell.Value = 45183
    }
}
